# Generate Report for Handback
# Refresh the handback status workbook with the results of a newer
# handback run: new source UUIDs, new xliff hash, and new timestamps.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "b84f5bbe-6c67-4ae1-b690-2f3d283370b5"
$oldUuid2 = "ee2b127f-e981-41b1-82e4-fb07c0b804a4"
$newUuid1 = "a0f10a7c-6e84-4123-8ad9-5bd465a4c7ad"
$newUuid2 = "ffffe662beed-9921-4904-9775-ac1f3056ca94"

$newHash = "d2780beeddf58ea9ef886415d987fd9b265baee5"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("A3").Value = "$newUuid2.md"

$wsOverview.Range("G2").Value = "2016-08-15 09:16:19"
$wsOverview.Range("G3").Value = "2016-08-15 09:16:19"

# Hyperlinks: this runtime's Hyperlinks.Delete() removes every hyperlink on
# the sheet, so rebuild the full set with updated display text but the
# same (unchanged) target addresses.
$wsOverview.Range("A1").Hyperlinks.Delete()

$wsOverview.Range("B2").Value = "e2e\$newUuid1.md"
$wsOverview.Range("B3").Value = "e2e\$newUuid2.md"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8a67c6b3b985181b261f687657507234b5c8a4b/e2e/$oldUuid1.md", "", "", "e2e\$newUuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8a67c6b3b985181b261f687657507234b5c8a4b/e2e/$oldUuid2.md", "", "", "e2e\$newUuid2.md")

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newUuid1.md"
$wsZhCn.Range("I2").Value = "$newUuid1.md"
$wsZhCn.Range("A3").Value = "$newUuid2.md"
$wsZhCn.Range("I3").Value = "$newUuid2.md"

$wsZhCn.Range("G2").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("J2").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-15 09:16:14"
$wsZhCn.Range("K2").Value = "2016-08-15 09:16:31"

$wsZhCn.Range("G3").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-15 09:16:14"
$wsZhCn.Range("K3").Value = "2016-08-15 09:16:31"

$wsZhCn.Range("A1").Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8a67c6b3b985181b261f687657507234b5c8a4b/e2e/$oldUuid1.md", "", "", "$newUuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b7b8830044e7759b719c39d57d18c788aa3a17c2/e2e/$oldUuid1.md", "", "", "$newUuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8a67c6b3b985181b261f687657507234b5c8a4b/e2e/$oldUuid2.md", "", "", "$newUuid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b7b8830044e7759b719c39d57d18c788aa3a17c2/e2e/$oldUuid2.md", "", "", "$newUuid2.md")

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newUuid1.md"
$wsDeDe.Range("I2").Value = "$newUuid1.md"
$wsDeDe.Range("A3").Value = "$newUuid2.md"
$wsDeDe.Range("I3").Value = "$newUuid2.md"

$wsDeDe.Range("G2").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("J2").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-15 09:16:19"
$wsDeDe.Range("K2").Value = "2016-08-15 09:16:38"

$wsDeDe.Range("G3").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("J3").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-15 09:16:19"
$wsDeDe.Range("K3").Value = "2016-08-15 09:16:38"

$wsDeDe.Range("A1").Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8a67c6b3b985181b261f687657507234b5c8a4b/e2e/$oldUuid1.md", "", "", "$newUuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/da3a48ba24ff612df0819ac3accbd4ffec479e50/e2e/$oldUuid1.md", "", "", "$newUuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8a67c6b3b985181b261f687657507234b5c8a4b/e2e/$oldUuid2.md", "", "", "$newUuid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/da3a48ba24ff612df0819ac3accbd4ffec479e50/e2e/$oldUuid2.md", "", "", "$newUuid2.md")

Write-Output "Handback status report regenerated."
